# Minijobbeitraege.xlsx - implement the "Arbeitnehmer Rentenpauschale" question
# for Minijobs and rename the flat-rate contribution labels accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row 3 ("Zahlt Arbeitnehmer Rentenpauschale?") - this pushes the
# former rows 3..10 down to 4..11 and keeps their number formatting (style)
# because Excel copies the formatting of the row above on insert.
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = "Zahlt Arbeitnehmer Rentenpauschale?"
$ws.Range("B3").Value = "ja"

# Rename the three flat-rate ("Pauschale") contribution labels that used to be
# called "...beitrag..." to reflect the new wording.
$ws.Range("A4").Value = "Arbeitgeberpauschale Krankenversicherung in Prozent"
$ws.Range("A5").Value = "Arbeitgeberpauschale Rentenversicherung in Prozent"
$ws.Range("A6").Value = "Arbeitnehmerpauschale Rentenversicherung in Prozent"

# Extend the "ja"/"nein" list validation (originally only on B2) to also cover
# the new B3 cell.
$null = $ws.Range("B2:B3").Validation.Add(3, 1, 1, "=Tabelle2!`$A`$2:`$A`$3")

# Update the active selection to match the author's final cursor position.
$null = $ws.Activate()
$null = $ws.Range("A8").Select()
